# Updates "horarios-141" schedule scrape workbook to the next scrape
# (Última actualización: 10:05:51), appending newly-scraped arrivals to
# each of the three sheets and re-sorting rows 79+ on sheet "LP1912" by
# arrival time (column B), same as the source scraper does on every run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "LP1912" (main sheet)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value = "Última actualización: 10:05:51"
$ws1.Cells.Item(3, 1).Value = "Total filas: 93"

# Two rows that were already present swap their "Linea" label (C40/C41)
$ws1.Cells.Item(40, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(41, 3).Value = "15_ABASTO"

# From row 79 onward, the table is re-sorted by Hora_Llegada after the
# new scrape is merged in; rewrite that whole tail with the final,
# sorted content (8 pre-existing rows + 12 freshly scraped ones).
$sheet1Tail = @(
    @("10:05:51", "10:43", "11X44_ETCHEVERRY",    "38",  "LP1912"),
    @("08:45:36", "10:44", "11X44_ETCHEVERRY",    "119", "LP1912"),
    @("08:52:50", "10:46", "15_P INDUSTRIAL",     "114", "LP1912"),
    @("10:05:51", "10:55", "16_SANTA ANA",        "50",  "LP1912"),
    @("10:05:51", "10:56", "27_EL RETIRO",        "51",  "LP1912"),
    @("09:23:23", "10:57", "10_OLMOS",            "94",  "LP1912"),
    @("10:05:51", "10:58", "10_OLMOS",            "53",  "LP1912"),
    @("09:23:23", "10:59", "27_EL RETIRO",        "96",  "LP1912"),
    @("09:23:23", "11:01", "81_EL PELIGRO",       "98",  "LP1912"),
    @("10:05:51", "11:04", "23_HERNANDEZ",        "59",  "LP1912"),
    @("09:23:23", "11:10", "16_P MOR-SANTA ANA",  "107", "LP1912"),
    @("09:23:23", "11:14", "14_ABASTO",           "111", "LP1912"),
    @("09:23:23", "11:15", "15X38_ABASTO",        "112", "LP1912"),
    @("10:05:51", "11:28", "10_OLMOS",            "83",  "LP1912"),
    @("10:05:51", "11:30", "215C_EL PATO",        "85",  "LP1912"),
    @("10:05:51", "11:31", "16_SANTA ANA",        "86",  "LP1912"),
    @("10:05:51", "11:41", "215B_EL PATO",        "96",  "LP1912"),
    @("10:05:51", "11:45", "15X38_ABASTO",        "100", "LP1912"),
    @("10:05:51", "11:52", "225_GOMEZ",           "107", "LP1912"),
    @("10:05:51", "11:58", "17_ROMERO",           "113", "LP1912")
)

$r = 79
foreach ($row in $sheet1Tail) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = [int]$row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value = "Última actualización: 10:05:51"
$ws2.Cells.Item(3, 1).Value = "Total filas: 19"

$ws2.Cells.Item(23, 1).Value = "10:05:51"
$ws2.Cells.Item(23, 2).Value = "11:30"
$ws2.Cells.Item(23, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(23, 4).Value = 85
$ws2.Cells.Item(23, 5).Value = "LP1912"

$ws2.Cells.Item(24, 1).Value = "10:05:51"
$ws2.Cells.Item(24, 2).Value = "11:41"
$ws2.Cells.Item(24, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(24, 4).Value = 96
$ws2.Cells.Item(24, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value = "Última actualización: 10:05:51"

$ws3.Cells.Item(22, 1).Value = "10:05:51"
$ws3.Cells.Item(22, 2).Value = "11:25"
$ws3.Cells.Item(22, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(22, 4).Value = 80
$ws3.Cells.Item(22, 5).Value = "L6203"
